$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'64.769.87"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.71%  "

$c = $ws.Range("D3")
$c.Value = "'3.145.11"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").Value = "  +0.00%  "

$c = $ws.Range("D5")
$c.Value = "'580.92"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "

$c = $ws.Range("D6")
$c.Value = "'147.20"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.64%  "

$ws.Range("E7").Value = "  +0.11%  "

$c = $ws.Range("D8")
$c.Value = "'3.143.65"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "

$c = $ws.Range("D9")
$c.Value = "'0.527"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "

$c = $ws.Range("D10")
$c.Value = "'0.158"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.45%  "

$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("E12").Value = "  -2.19%  "

$ws.Range("E13").Value = "  -3.09%  "

$ws.Range("E14").Value = "  -3.06%  "

$c = $ws.Range("D15")
$c.Value = "'3.660.03"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "

$c = $ws.Range("D16")
$c.Value = "'64.839.34"
$c.Style = "Normal"

$ws.Range("E17").Value = "  -1.33%  "

$c = $ws.Range("D18")
$c.Value = "'3.136.51"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("E19").Value = "  -0.54%  "

$c = $ws.Range("D20")
$c.Value = "'499.86"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.89%  "

$c = $ws.Range("D21")
$c.Value = "'15.44"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.14%  "

$ws.Range("E22").Value = "  -3.65%  "

$c = $ws.Range("D23")
$c.Value = "'14.98"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -7.18%  "

$c = $ws.Range("D24")
$c.Value = "'7.76"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.58%  "

$c = $ws.Range("D25")
$c.Value = "'84.31"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("E26").Value = "  +0.20%  "

$c = $ws.Range("D27")
$c.Value = "'9.05"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("E28").Value = "  +0.33%  "

$c = $ws.Range("D29")
$c.Value = "'2.18"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.14%  "

$ws.Range("E30").Value = "  +0.82%  "

$c = $ws.Range("D31")
$c.Value = "'27.58"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.84%  "

$ws.Range("E32").Value = "  -0.74%  "

$c = $ws.Range("D33")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "

$c = $ws.Range("D34")
$c.Value = "'6.38"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.20%  "

$c = $ws.Range("D35")
$c.Value = "'6.44"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.24%  "

$c = $ws.Range("D36")
$c.Value = "'54.82"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.01%  "

$c = $ws.Range("D37")
$c.Value = "'0.0892"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "

$c = $ws.Range("D38")
$c.Value = "'469.43"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("E39").Value = "  -1.14%  "

$c = $ws.Range("D40")
$c.Value = "'2.90"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -6.30%  "

$c = $ws.Range("D41")
$c.Value = "'8.73"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "

$c = $ws.Range("D42")
$c.Value = "'2.979.65"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.59%  "

$ws.Range("E43").Value = "  -3.88%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D44")
$c.Value = "'0.282"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.32%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D45")
$c.Value = "'2.42"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.92%  "

$c = $ws.Range("D46")
$c.Value = "'28.22"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.75%  "

$c = $ws.Range("D47")
$c.Value = "'0.0₃0599"
$c.Style = "Normal"

$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("E50").Value = "  -4.16%  "

$c = $ws.Range("D51")
$c.Value = "'119.00"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.29%  "
